# "this is commit for today"
#
# The original single paragraph reads (across 4 runs, with spell/grammar
# proofErr markers around "gansh" and " ."):
#   Hi this is |gansh| .  I also updated here
#
# It becomes a single merged run reading:
#   Hi this is gansh .  I also updated here
# and that paragraph is then repeated, separated by blank paragraphs, so
# the document ends up as 7 paragraphs: text, text, blank, text, blank,
# text, blank.

$d = $word.ActiveDocument

$targetText = "Hi this is gansh .  I also updated here "

# Step 1: rewrite the first (only) paragraph's text. Matching the whole
# original run-spanning text in one Find/Replace merges the runs (and
# drops the spellStart/gramStart/spellEnd/gramEnd proofErr markers)
# into a single run, while leaving the paragraph's own properties
# (w14:paraId, rsids, pPr, ...) untouched.
$d.Content.Find.Execute(
    "Hi this is gansh .  I also updated here ",
    $true, $false, $false, $false, $false, $true, 1, $false,
    $targetText, 2) | Out-Null

# Step 2: append the 6 remaining paragraphs (text / blank / text / blank /
# text / blank) as raw OOXML right after the first paragraph. Building
# these by hand (rather than via InsertParagraphAfter, which leaves a
# stray empty <w:r> behind in blank paragraphs) keeps the blank
# paragraphs truly empty, matching the target markup exactly.
$wNs = "xmlns:w='http://schemas.openxmlformats.org/wordprocessingml/2006/main'"
$pWithText = "<w:p $wNs><w:pPr><w:rPr><w:lang w:val='en-US'/></w:rPr></w:pPr><w:r><w:rPr><w:lang w:val='en-US'/></w:rPr><w:t xml:space='preserve'>$targetText</w:t></w:r></w:p>"
$pBlank = "<w:p $wNs><w:pPr><w:rPr><w:lang w:val='en-US'/></w:rPr></w:pPr></w:p>"

$fragment = $pWithText + $pBlank + $pWithText + $pBlank + $pWithText + $pBlank

$r = $d.Content
$r.Collapse(0) | Out-Null
$r.InsertXML($fragment)

Write-Output ("paragraphs: " + $d.Paragraphs.Count)
Write-Output $d.Content.Text
